$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 5.795774333333334
$ws.Cells.Item(2, 8).Value = 17.387323
$ws.Cells.Item(2, 9).Value = 0.2166180251727471
$ws.Cells.Item(2, 10).Value = 0.2166180251727471
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 1.181145333333333
$ws.Cells.Item(2, 14).Value = 3.543436
$ws.Cells.Item(2, 15).Value = 0.1469875775727626
$ws.Cells.Item(2, 16).Value = 0.1469875775727626
$ws.Cells.Item(2, 17).Value = 6.845651806869779
$ws.Cells.Item(2, 18).Value = 61.61086626182801
$ws.Cells.Item(2, 19).Value = 0.0318401587787378
$ws.Cells.Item(2, 20).Value = 0.0318401587787378

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 5.795774333333334
$ws.Cells.Item(3, 8).Value = 17.387323
$ws.Cells.Item(3, 9).Value = 0.2166180251727471
$ws.Cells.Item(3, 10).Value = 0.2166180251727471
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 5.122951333333333
$ws.Cells.Item(3, 14).Value = 15.368854
$ws.Cells.Item(3, 15).Value = 0.637525446919166
$ws.Cells.Item(3, 16).Value = 0.637525446919166
$ws.Cells.Item(3, 17).Value = 29.69146984864911
$ws.Cells.Item(3, 18).Value = 267.223228637842
$ws.Cells.Item(3, 19).Value = 0.1380995033090028
$ws.Cells.Item(3, 20).Value = 0.1380995033090028

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 5.795774333333334
$ws.Cells.Item(4, 8).Value = 17.387323
$ws.Cells.Item(4, 9).Value = 0.2166180251727471
$ws.Cells.Item(4, 10).Value = 0.2166180251727471
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 1.731584666666667
$ws.Cells.Item(4, 14).Value = 5.194754
$ws.Cells.Item(4, 15).Value = 0.2154869755080714
$ws.Cells.Item(4, 16).Value = 0.2154869755080714
$ws.Cells.Item(4, 17).Value = 10.03587396706022
$ws.Cells.Item(4, 18).Value = 90.32286570354201
$ws.Cells.Item(4, 19).Value = 0.04667836308500656
$ws.Cells.Item(4, 20).Value = 0.04667836308500656

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 4.872962999999999
$ws.Cells.Item(5, 8).Value = 14.618889
$ws.Cells.Item(5, 9).Value = 0.1821277988221416
$ws.Cells.Item(5, 10).Value = 0.1821277988221416
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 1.181145333333333
$ws.Cells.Item(5, 14).Value = 3.543436
$ws.Cells.Item(5, 15).Value = 0.1469875775727626
$ws.Cells.Item(5, 16).Value = 0.1469875775727626
$ws.Cells.Item(5, 17).Value = 5.755677506956
$ws.Cells.Item(5, 18).Value = 51.801097562604
$ws.Cells.Item(5, 19).Value = 0.02677052395752603
$ws.Cells.Item(5, 20).Value = 0.02677052395752603

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 4.872962999999999
$ws.Cells.Item(6, 8).Value = 14.618889
$ws.Cells.Item(6, 9).Value = 0.1821277988221416
$ws.Cells.Item(6, 10).Value = 0.1821277988221416
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 5.122951333333333
$ws.Cells.Item(6, 14).Value = 15.368854
$ws.Cells.Item(6, 15).Value = 0.637525446919166
$ws.Cells.Item(6, 16).Value = 0.637525446919166
$ws.Cells.Item(6, 17).Value = 24.963952298134
$ws.Cells.Item(6, 18).Value = 224.675570683206
$ws.Cells.Item(6, 19).Value = 0.1161111063404898
$ws.Cells.Item(6, 20).Value = 0.1161111063404898

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 4.872962999999999
$ws.Cells.Item(7, 8).Value = 14.618889
$ws.Cells.Item(7, 9).Value = 0.1821277988221416
$ws.Cells.Item(7, 10).Value = 0.1821277988221416
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 1.731584666666667
$ws.Cells.Item(7, 14).Value = 5.194754
$ws.Cells.Item(7, 15).Value = 0.2154869755080714
$ws.Cells.Item(7, 16).Value = 0.2154869755080714
$ws.Cells.Item(7, 17).Value = 8.437948012033999
$ws.Cells.Item(7, 18).Value = 75.94153210830599
$ws.Cells.Item(7, 19).Value = 0.03924616852412578
$ws.Cells.Item(7, 20).Value = 0.03924616852412578

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 4.523564
$ws.Cells.Item(8, 8).Value = 13.570692
$ws.Cells.Item(8, 9).Value = 0.1690689533557062
$ws.Cells.Item(8, 10).Value = 0.1690689533557062
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 1.181145333333333
$ws.Cells.Item(8, 14).Value = 3.543436
$ws.Cells.Item(8, 15).Value = 0.1469875775727626
$ws.Cells.Item(8, 16).Value = 0.1469875775727626
$ws.Cells.Item(8, 17).Value = 5.342986508634668
$ws.Cells.Item(8, 18).Value = 48.08687857771201
$ws.Cells.Item(8, 19).Value = 0.02485103589651764
$ws.Cells.Item(8, 20).Value = 0.02485103589651764

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 4.523564
$ws.Cells.Item(9, 8).Value = 13.570692
$ws.Cells.Item(9, 9).Value = 0.1690689533557062
$ws.Cells.Item(9, 10).Value = 0.1690689533557062
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 5.122951333333333
$ws.Cells.Item(9, 14).Value = 15.368854
$ws.Cells.Item(9, 15).Value = 0.637525446919166
$ws.Cells.Item(9, 16).Value = 0.637525446919166
$ws.Cells.Item(9, 17).Value = 23.17399822521867
$ws.Cells.Item(9, 18).Value = 208.565984026968
$ws.Cells.Item(9, 19).Value = 0.1077857600482522
$ws.Cells.Item(9, 20).Value = 0.1077857600482522

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 4.523564
$ws.Cells.Item(10, 8).Value = 13.570692
$ws.Cells.Item(10, 9).Value = 0.1690689533557062
$ws.Cells.Item(10, 10).Value = 0.1690689533557062
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 1.731584666666667
$ws.Cells.Item(10, 14).Value = 5.194754
$ws.Cells.Item(10, 15).Value = 0.2154869755080714
$ws.Cells.Item(10, 16).Value = 0.2154869755080714
$ws.Cells.Item(10, 17).Value = 7.832934061085333
$ws.Cells.Item(10, 18).Value = 70.496406549768
$ws.Cells.Item(10, 19).Value = 0.03643215741093633
$ws.Cells.Item(10, 20).Value = 0.03643215741093633

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 5.027470333333333
$ws.Cells.Item(11, 8).Value = 15.082411
$ws.Cells.Item(11, 9).Value = 0.1879025359834701
$ws.Cells.Item(11, 10).Value = 0.1879025359834701
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 1.181145333333333
$ws.Cells.Item(11, 14).Value = 3.543436
$ws.Cells.Item(11, 15).Value = 0.1469875775727626
$ws.Cells.Item(11, 16).Value = 0.1469875775727626
$ws.Cells.Item(11, 17).Value = 5.938173122688445
$ws.Cells.Item(11, 18).Value = 53.44355810419601
$ws.Cells.Item(11, 19).Value = 0.02761933858398912
$ws.Cells.Item(11, 20).Value = 0.02761933858398912

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 5.027470333333333
$ws.Cells.Item(12, 8).Value = 15.082411
$ws.Cells.Item(12, 9).Value = 0.1879025359834701
$ws.Cells.Item(12, 10).Value = 0.1879025359834701
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 5.122951333333333
$ws.Cells.Item(12, 14).Value = 15.368854
$ws.Cells.Item(12, 15).Value = 0.637525446919166
$ws.Cells.Item(12, 16).Value = 0.637525446919166
$ws.Cells.Item(12, 17).Value = 25.75548584744378
$ws.Cells.Item(12, 18).Value = 231.799372626994
$ws.Cells.Item(12, 19).Value = 0.1197926482301065
$ws.Cells.Item(12, 20).Value = 0.1197926482301065

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 5.027470333333333
$ws.Cells.Item(13, 8).Value = 15.082411
$ws.Cells.Item(13, 9).Value = 0.1879025359834701
$ws.Cells.Item(13, 10).Value = 0.1879025359834701
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 1.731584666666667
$ws.Cells.Item(13, 14).Value = 5.194754
$ws.Cells.Item(13, 15).Value = 0.2154869755080714
$ws.Cells.Item(13, 16).Value = 0.2154869755080714
$ws.Cells.Item(13, 17).Value = 8.705490541321556
$ws.Cells.Item(13, 18).Value = 78.349414871894
$ws.Cells.Item(13, 19).Value = 0.04049054916937454
$ws.Cells.Item(13, 20).Value = 0.04049054916937454

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 3.600244
$ws.Cells.Item(14, 8).Value = 10.800732
$ws.Cells.Item(14, 9).Value = 0.1345597155042265
$ws.Cells.Item(14, 10).Value = 0.1345597155042266
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 1.181145333333333
$ws.Cells.Item(14, 14).Value = 3.543436
$ws.Cells.Item(14, 15).Value = 0.1469875775727626
$ws.Cells.Item(14, 16).Value = 0.1469875775727626
$ws.Cells.Item(14, 17).Value = 4.252411399461334
$ws.Cells.Item(14, 18).Value = 38.271702595152
$ws.Cells.Item(14, 19).Value = 0.01977860662084636
$ws.Cells.Item(14, 20).Value = 0.01977860662084636

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 3.600244
$ws.Cells.Item(15, 8).Value = 10.800732
$ws.Cells.Item(15, 9).Value = 0.1345597155042265
$ws.Cells.Item(15, 10).Value = 0.1345597155042266
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 5.122951333333333
$ws.Cells.Item(15, 14).Value = 15.368854
$ws.Cells.Item(15, 15).Value = 0.637525446919166
$ws.Cells.Item(15, 16).Value = 0.637525446919166
$ws.Cells.Item(15, 17).Value = 18.44387480012533
$ws.Cells.Item(15, 18).Value = 165.994873201128
$ws.Cells.Item(15, 19).Value = 0.08578524276414787
$ws.Cells.Item(15, 20).Value = 0.08578524276414788

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 3.600244
$ws.Cells.Item(16, 8).Value = 10.800732
$ws.Cells.Item(16, 9).Value = 0.1345597155042265
$ws.Cells.Item(16, 10).Value = 0.1345597155042266
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 1.731584666666667
$ws.Cells.Item(16, 14).Value = 5.194754
$ws.Cells.Item(16, 15).Value = 0.2154869755080714
$ws.Cells.Item(16, 16).Value = 0.2154869755080714
$ws.Cells.Item(16, 17).Value = 6.234127306658666
$ws.Cells.Item(16, 18).Value = 56.107145759928
$ws.Cells.Item(16, 19).Value = 0.02899586611923232
$ws.Cells.Item(16, 20).Value = 0.02899586611923233

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 2.935718666666666
$ws.Cells.Item(17, 8).Value = 8.807155999999999
$ws.Cells.Item(17, 9).Value = 0.1097229711617085
$ws.Cells.Item(17, 10).Value = 0.1097229711617085
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 1.181145333333333
$ws.Cells.Item(17, 14).Value = 3.543436
$ws.Cells.Item(17, 15).Value = 0.1469875775727626
$ws.Cells.Item(17, 16).Value = 0.1469875775727626
$ws.Cells.Item(17, 17).Value = 3.467510403112889
$ws.Cells.Item(17, 18).Value = 31.207593628016
$ws.Cells.Item(17, 19).Value = 0.01612791373514562
$ws.Cells.Item(17, 20).Value = 0.01612791373514562

$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 2.935718666666666
$ws.Cells.Item(18, 8).Value = 8.807155999999999
$ws.Cells.Item(18, 9).Value = 0.1097229711617085
$ws.Cells.Item(18, 10).Value = 0.1097229711617085
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 5.122951333333333
$ws.Cells.Item(18, 14).Value = 15.368854
$ws.Cells.Item(18, 15).Value = 0.637525446919166
$ws.Cells.Item(18, 16).Value = 0.637525446919166
$ws.Cells.Item(18, 17).Value = 15.03954385769155
$ws.Cells.Item(18, 18).Value = 135.355894719224
$ws.Cells.Item(18, 19).Value = 0.06995118622716695
$ws.Cells.Item(18, 20).Value = 0.06995118622716695

$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 2.935718666666666
$ws.Cells.Item(19, 8).Value = 8.807155999999999
$ws.Cells.Item(19, 9).Value = 0.1097229711617085
$ws.Cells.Item(19, 10).Value = 0.1097229711617085
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 13).Value = 1.731584666666667
$ws.Cells.Item(19, 14).Value = 5.194754
$ws.Cells.Item(19, 15).Value = 0.2154869755080714
$ws.Cells.Item(19, 16).Value = 0.2154869755080714
$ws.Cells.Item(19, 17).Value = 5.08344542884711
$ws.Cells.Item(19, 18).Value = 45.75100885962399
$ws.Cells.Item(19, 19).Value = 0.0236438711993959
$ws.Cells.Item(19, 20).Value = 0.0236438711993959
